$d = $word.ActiveDocument

function Protect-Boundary($rng) {
    # Force the given range to remain (or become) its own run, distinct
    # from neighbouring runs, by toggling a character format on then off.
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1. Title: "QEDistribution Widget" -> "QENumericEdit Widget"
#    Keep the surrounding runs ("QE", " ", "Widget") intact/unmerged.
# ---------------------------------------------------------------------
$titleWord = $d.Range(3, 15)
$titleWord.Text = "NumericEdit"
Protect-Boundary $d.Range(3, 14)
Protect-Boundary $d.Range(14, 15)

# ---------------------------------------------------------------------
# 2. Date line: "25th July 2020" -> "16th July 2020"
#    Runs "2" + "5" collapse into a single "16" run; "July"," ","20","20"
#    must remain four separate runs as before.
# ---------------------------------------------------------------------
Protect-Boundary $d.Range(44, 48)   # "July"
Protect-Boundary $d.Range(48, 49)   # " "
Protect-Boundary $d.Range(49, 51)   # "20"
Protect-Boundary $d.Range(51, 53)   # "20"
$dateRange = $d.Range(39, 41)
$dateRange.Text = "16"

Write-Output $d.Paragraphs(2).Range.Text
Write-Output $d.Paragraphs(4).Range.Text
